$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "0+12=12"
$t.Cell(1,2).Range.Text = "94-13=81"
$t.Cell(1,3).Range.Text = "21+36=57"
$t.Cell(1,4).Range.Text = "69-4=65"
$t.Cell(1,5).Range.Text = "44-34=10"
$t.Cell(2,1).Range.Text = "38-23=15"
$t.Cell(2,2).Range.Text = "14-10=4"
$t.Cell(2,3).Range.Text = "94-89=5"
$t.Cell(2,4).Range.Text = "80+13=93"
$t.Cell(2,5).Range.Text = "66-29=37"
$t.Cell(3,1).Range.Text = "72-29=43"
$t.Cell(3,2).Range.Text = "34-18=16"
$t.Cell(3,3).Range.Text = "39+36=75"
$t.Cell(3,4).Range.Text = "24+15=39"
$t.Cell(3,5).Range.Text = "92-1=91"
$t.Cell(4,1).Range.Text = "25+30=55"
$t.Cell(4,2).Range.Text = "54+35=89"
$t.Cell(4,3).Range.Text = "1+65=66"
$t.Cell(4,4).Range.Text = "25+64=89"
$t.Cell(4,5).Range.Text = "1+94=95"
$t.Cell(5,1).Range.Text = "74-5=69"
$t.Cell(5,2).Range.Text = "81-40=41"
$t.Cell(5,3).Range.Text = "10+38=48"
$t.Cell(5,4).Range.Text = "20+4=24"
$t.Cell(5,5).Range.Text = "22+28=50"
$t.Cell(6,1).Range.Text = "91-66=25"
$t.Cell(6,2).Range.Text = "22+49=71"
$t.Cell(6,3).Range.Text = "25-22=3"
$t.Cell(6,4).Range.Text = "9+35=44"
$t.Cell(6,5).Range.Text = "36+51=87"
$t.Cell(7,1).Range.Text = "37+39=76"
$t.Cell(7,2).Range.Text = "59+32=91"
$t.Cell(7,3).Range.Text = "37-27=10"
$t.Cell(7,4).Range.Text = "91+6=97"
$t.Cell(7,5).Range.Text = "74+20=94"
$t.Cell(8,1).Range.Text = "64+32=96"
$t.Cell(8,2).Range.Text = "36-26=10"
$t.Cell(8,3).Range.Text = "83+5=88"
$t.Cell(8,4).Range.Text = "68-51=17"
$t.Cell(8,5).Range.Text = "69-52=17"
$t.Cell(9,1).Range.Text = "81-3=78"
$t.Cell(9,2).Range.Text = "16+0=16"
$t.Cell(9,3).Range.Text = "34-24=10"
$t.Cell(9,4).Range.Text = "51+0=51"
$t.Cell(9,5).Range.Text = "38+18=56"
$t.Cell(10,1).Range.Text = "96+0=96"
$t.Cell(10,2).Range.Text = "28+69=97"
$t.Cell(10,3).Range.Text = "37+26=63"
$t.Cell(10,4).Range.Text = "26-6=20"
$t.Cell(10,5).Range.Text = "44+32=76"
$t.Cell(11,1).Range.Text = "24+63=87"
$t.Cell(11,2).Range.Text = "95-88=7"
$t.Cell(11,3).Range.Text = "60-53=7"
$t.Cell(11,4).Range.Text = "31-22=9"
$t.Cell(11,5).Range.Text = "29+59=88"
$t.Cell(12,1).Range.Text = "9+47=56"
$t.Cell(12,2).Range.Text = "41-28=13"
$t.Cell(12,3).Range.Text = "62-1=61"
$t.Cell(12,4).Range.Text = "51+27=78"
$t.Cell(12,5).Range.Text = "46+22=68"
$t.Cell(13,1).Range.Text = "30+41=71"
$t.Cell(13,2).Range.Text = "66-1=65"
$t.Cell(13,3).Range.Text = "23+73=96"
$t.Cell(13,4).Range.Text = "92-45=47"
$t.Cell(13,5).Range.Text = "20-11=9"
$t.Cell(14,1).Range.Text = "87-48=39"
$t.Cell(14,2).Range.Text = "12+87=99"
$t.Cell(14,3).Range.Text = "0+59=59"
$t.Cell(14,4).Range.Text = "86-78=8"
$t.Cell(14,5).Range.Text = "95-90=5"
$t.Cell(15,1).Range.Text = "13+22=35"
$t.Cell(15,2).Range.Text = "33+38=71"
$t.Cell(15,3).Range.Text = "24+58=82"
$t.Cell(15,4).Range.Text = "76-7=69"
$t.Cell(15,5).Range.Text = "80-13=67"
$t.Cell(16,1).Range.Text = "70-62=8"
$t.Cell(16,2).Range.Text = "70-6=64"
$t.Cell(16,3).Range.Text = "23-13=10"
$t.Cell(16,4).Range.Text = "28+2=30"
$t.Cell(16,5).Range.Text = "11+57=68"
$t.Cell(17,1).Range.Text = "9+73=82"
$t.Cell(17,2).Range.Text = "23+64=87"
$t.Cell(17,3).Range.Text = "77+10=87"
$t.Cell(17,4).Range.Text = "19+63=82"
$t.Cell(17,5).Range.Text = "9+69=78"
$t.Cell(18,1).Range.Text = "33+2=35"
$t.Cell(18,2).Range.Text = "91-35=56"
$t.Cell(18,3).Range.Text = "61-24=37"
$t.Cell(18,4).Range.Text = "69-26=43"
$t.Cell(18,5).Range.Text = "74-32=42"
$t.Cell(19,1).Range.Text = "8+47=55"
$t.Cell(19,2).Range.Text = "11+10=21"
$t.Cell(19,3).Range.Text = "28-7=21"
$t.Cell(19,4).Range.Text = "81+3=84"
$t.Cell(19,5).Range.Text = "21+21=42"
$t.Cell(20,1).Range.Text = "98-74=24"
$t.Cell(20,2).Range.Text = "22+50=72"
$t.Cell(20,3).Range.Text = "23-23=0"
$t.Cell(20,4).Range.Text = "92-70=22"
$t.Cell(20,5).Range.Text = "64+6=70"
